# ooSQLite - make some arguments optional; otherwise changes are formatting changes
#
# sqlite3_backup_init / sqlite3_backup_pagecount / sqlite3_backup_remaining /
# sqlite3_backup_step were previously marked "planned" with a comment saying
# "need to implement ooSQLiteBackup class". They are now marked as
# "implemented oo and classic" (same status/style as the row right above,
# sqlite3_backup_finish), and the now-obsolete comment is cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 9, 10, 11, 12

foreach ($r in $rows) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    # Copy the "implemented" look (font/fill) from B8, which already has it.
    $ws.Range("B8").Copy()
    $bCell.PasteSpecial(-4122)  # xlPasteFormats

    # Update the status text and clear out the no-longer-needed comment.
    $bCell.Value = " implemented oo and classic"
    $cCell.Value = ""
}

$excel.CutCopyMode = $false

# Restore the view: scrolled down with C9:C12 selected, matching the edited rows.
$ws.Range("C9:C12").Select()
$excel.ActiveWindow.ScrollRow = 185
